# Update the "江西-漫展信息" workbook:
#  - bump a handful of "想去人数" (F column) counters on both the
#    "展览" and "全部类型" sheets (they carry duplicate data tables);
#  - insert a new Nanchang "CM03 配音演员孙路路专场见面会" event as the
#    new row 36, pushing the three rows that used to sit at 36-38 down
#    to 37-39 (the running index in column A is a plain serial and is
#    NOT renumbered - it still reads 35/36/37 for those three rows, and
#    the freshly appended row 39 gets 38).

function Set-TextValue($range, $text) {
    # Several columns hold strings that look like dates / date ranges
    # (e.g. "2024-08-17", "2024.08.17 09:00-08.17 17:00"). A bare
    # Range.Value assignment lets Excel's type-inference turn those into
    # real date serials, so force a text format for the write and then
    # drop the formatting again so the cell ends up styled exactly like
    # its neighbours (no explicit numFmt left behind).
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

function Update-ConManSheet($ws) {
    # ---- 1. small "想去人数" (column F) bumps across existing rows ----
    $fUpdates = @{
        3  = 5040
        5  = 7314
        8  = 98
        12 = 4266
        13 = 1722
        15 = 95
        16 = 2871
        19 = 203
        21 = 422
        22 = 447
        23 = 281
        24 = 86
        25 = 1668
        27 = 85
        28 = 1353
        34 = 55
        35 = 102
    }
    foreach ($r in $fUpdates.Keys) {
        $ws.Cells.Item($r, 6).Value = $fUpdates[$r]
    }

    # ---- 2. capture the current rows 36-38 (B:I) before they move ----
    $oldRows = @(36, 37, 38)
    $captured = @{}
    foreach ($r in $oldRows) {
        $captured[$r] = @(
            $ws.Cells.Item($r, 2).Value(),
            $ws.Cells.Item($r, 3).Value(),
            $ws.Cells.Item($r, 4).Value(),
            $ws.Cells.Item($r, 5).Value(),
            $ws.Cells.Item($r, 6).Value(),
            $ws.Cells.Item($r, 7).Value(),
            $ws.Cells.Item($r, 8).Value(),
            $ws.Cells.Item($r, 9).Value()
        )
    }

    # ---- 3. give row 39 the same look as the other index cells ----
    $ws.Range("A38").Copy()
    $ws.Range("A39").PasteSpecial(-4122)   # xlPasteFormats
    $ws.Range("A39").Value = 38

    # ---- 4. write the captured rows one row lower than they used to be ----
    # old row 36 (CM03 展览会)      -> new row 37, with F updated 2713->2736
    # old row 37 (第四届龙年动漫展) -> new row 38, with F updated 690->694
    # old row 38 (哥布林动漫游戏展) -> new row 39, F unchanged (37)
    $shiftMap = @{ 36 = 37; 37 = 38; 38 = 39 }
    $fOverride = @{ 37 = 2736; 38 = 694 }

    foreach ($oldR in @(38, 37, 36)) {
        $newR = $shiftMap[$oldR]
        $vals = $captured[$oldR]

        Set-TextValue $ws.Cells.Item($newR, 2) $vals[0]
        $ws.Cells.Item($newR, 3).Value = $vals[1]
        $ws.Cells.Item($newR, 4).Value = $vals[2]
        Set-TextValue $ws.Cells.Item($newR, 5) $vals[3]

        if ($fOverride.ContainsKey($newR)) {
            $ws.Cells.Item($newR, 6).Value = $fOverride[$newR]
        } else {
            $ws.Cells.Item($newR, 6).Value = $vals[4]
        }
        $ws.Cells.Item($newR, 7).Value = $vals[5]
        $ws.Cells.Item($newR, 8).Value = $vals[6]
        $ws.Cells.Item($newR, 9).Value = $vals[7]
    }

    # ---- 5. fill the now-vacant row 36 with the new event ----
    Set-TextValue $ws.Cells.Item(36, 2) "2024-08-17"
    $ws.Cells.Item(36, 3).Value = "南昌·CM03·配音演员孙路路专场见面会"
    $ws.Cells.Item(36, 4).Value = "怀玉山大道1315号 南昌绿地国际博览中心"
    Set-TextValue $ws.Cells.Item(36, 5) "2024.08.17 09:00-08.17 17:00"
    $ws.Cells.Item(36, 6).Value = 23
    $ws.Cells.Item(36, 7).Value = 258
    $ws.Cells.Item(36, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89623"
    $ws.Cells.Item(36, 9).Value = "//i2.hdslb.com/bfs/openplatform/202407/TIDNWGni1721565446649.jpeg"
}

$wb = $excel.ActiveWorkbook

Update-ConManSheet $wb.Worksheets.Item("展览")
Update-ConManSheet $wb.Worksheets.Item("全部类型")
